# Parameterization is added for all testcases
# Fills in previously-blank parameter cells on a few of the test-data
# sheets (First/Last Name + Employee ID, UserName, etc.) and switches the
# active tab from RecruitmentCandidate to EditEmployee.

$wb = $excel.ActiveWorkbook

# --- AddEmployee sheet: add First Name / Last Name / Employee ID for rows 2 & 3 ---
$wsAddEmployee = $wb.Worksheets.Item("AddEmployee")
$wsAddEmployee.Range("F2").Value = "Val"
$wsAddEmployee.Range("G2").Value = "Johnson"
$wsAddEmployee.Range("M2").Value = "0150"

$wsAddEmployee.Range("F3").Value = "Heriberto"
$wsAddEmployee.Range("G3").Value = "Moore"
$wsAddEmployee.Range("M3").Value = "0149"

# --- AddUser sheet: add Employee ID / UserName for row 3 ---
$wsAddUser = $wb.Worksheets.Item("AddUser")
$wsAddUser.Range("F3").Value = "0141"
$wsAddUser.Range("G3").Value = "herma.rodriguez"

# --- EditEmployee sheet: swap Employee ID / NewLastName / Location for row 6 ---
$wsEditEmployee = $wb.Worksheets.Item("EditEmployee")
$wsEditEmployee.Range("F6").Value = "1069"
$wsEditEmployee.Range("G6").Value = "Gibson"
$wsEditEmployee.Range("I6").Value = "Singapore Regional HQ"

# --- RecruitmentCandidate sheet: swap FirstName / LastName / LinkedInURL for row 6 ---
$wsRecruitmentCandidate = $wb.Worksheets.Item("RecruitmentCandidate")
$wsRecruitmentCandidate.Range("F6").Value = "Houston"
$wsRecruitmentCandidate.Range("G6").Value = "Hagenes"
$wsRecruitmentCandidate.Range("H6").Value = "Layla"

# --- Active tab moves from RecruitmentCandidate to EditEmployee, with a new selection ---
$wsEditEmployee.Activate()
$wsEditEmployee.Range("I6").Select()
